# "Building the process - Part Two"
# The workbook currently has two sheets:
#   - "Sheet1"   (sheetId 2, rId1) -> empty placeholder sheet
#   - "NewSheet" (sheetId 1, rId2) -> the real Coffee Shop report data
#
# The edit removes the empty placeholder "Sheet1" and renames the
# data-bearing "NewSheet" to "Sheet1", making it the single, active sheet
# in the workbook.

$wb = $excel.ActiveWorkbook

# Avoid the "delete sheet" confirmation prompt that real Excel would show.
$excel.DisplayAlerts = $false

# Remove the empty placeholder sheet.
$wb.Worksheets.Item("Sheet1").Delete()

# Rename the remaining (data) sheet to "Sheet1" and make it the active tab.
$dataSheet = $wb.Worksheets.Item("NewSheet")
$dataSheet.Name = "Sheet1"
$dataSheet.Activate()
$dataSheet.Select()

# Keep the selected cell on the data sheet consistent with the saved file
# (selection stays on C2).
$dataSheet.Range("C2").Select()

$excel.DisplayAlerts = $true
